# Apply updated Keyword/Correlation values to the Correlation Table SP500 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Keyword, Correlation)
$updates = @{
    2  = @("credit",  -0.6546)
    3  = @("trade",   -0.6784)
    4  = @("bank",    -0.7616000000000001)
    7  = @("bank",    -0.4116)
    8  = @("credit",  -0.589)
    9  = @("inflation", -0.4892)
    10 = @("trade",   0.1004)
    12 = @("credit",  -0.5551)
    13 = @("trade",   0.8056)
    14 = @("bank",    -0.3615)
    17 = @("bank",    -0.1084)
    18 = @("trade",   -0.6375999999999999)
    19 = @("inflation", 0.0639)
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
